# #3473 replaced two properties that had gaps
#
# Updates a handful of cells in the "BPS Data" worksheet:
#   - Row 2  (Medstar POB North Tower)  -> replaced with "Medstar POB South Tower" data
#   - Row 3  Owner name tweak ("1801 Pennsylvania Ave." -> "1801 Pennsylvania Avenue, LLC")
#   - Row 4  Address/Owner/Gross area corrections (GSA: 300 E Street SW)
#   - Row 5  Gross area correction (Paul H.Nitze)
#   - Row 6  (President Madison Apartments) -> replaced with "Hampton House" data
#   - Row 7  Postal code + gross area correction (3303 Water Street)
#   - Row 8  Address correction (15th and H Street Associates LLP)
#   - Row 9  unchanged (Eastern Market)
#   - Row 10 (DPW Vehicle Maintenance Facility 2) -> replaced with
#            "School Without Walls @ Francis Stevens" data
#
# Also clears the stray date number-format that had been applied to the
# "Year Built" column (I2:I10) -- it should just hold plain years.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2: Medstar POB North Tower -> Medstar POB South Tower
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = "Medstar POB South Tower"
$ws.Range("I2").Value = 1985
$ws.Range("L2").Value = 76319
$ws.Range("M2").Value = 12.8
$ws.Range("N2").Value = 150.2

# ---------------------------------------------------------------------------
# Row 3: owner name correction
# ---------------------------------------------------------------------------
$ws.Range("C3").Value = "1801 Pennsylvania Avenue, LLC"

# ---------------------------------------------------------------------------
# Row 4: GSA: 300 E Street SW
# ---------------------------------------------------------------------------
$ws.Range("E4").Value = "300 E ST SW"
$ws.Range("J4").Value = "TWO INDEPENDENCE HANA OW LLC"
$ws.Range("L4").Value = 627655

# ---------------------------------------------------------------------------
# Row 5: Paul H.Nitze - gross area correction
# ---------------------------------------------------------------------------
$ws.Range("L5").Value = 58717

# ---------------------------------------------------------------------------
# Row 6: President Madison Apartments -> Hampton House
# ---------------------------------------------------------------------------
$ws.Range("C6").Value = "Hampton House"
$ws.Range("E6").Value = "2700 CONNECTICUT AVENUE NW"
$ws.Range("H6").Value = 20008
$ws.Range("I6").Value = 1921
$ws.Range("J6").Value = "2700 CONECTICUT AVENUE LLC"
$ws.Range("L6").Value = 83580
$ws.Range("N6").Value = 59.1

# ---------------------------------------------------------------------------
# Row 7: 3303 Water Street - postal code + gross area correction
# ---------------------------------------------------------------------------
$ws.Range("H7").Value = 20005
$ws.Range("L7").Value = 145697

# ---------------------------------------------------------------------------
# Row 8: 15th and H Street Associates LLP - address correction
# ---------------------------------------------------------------------------
$ws.Range("E8").Value = "1428 H ST NW"

# ---------------------------------------------------------------------------
# Row 9: Eastern Market - no change
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Row 10: DPW Vehicle Maintenance Facility 2 -> School Without Walls @ Francis Stevens
# ---------------------------------------------------------------------------
$ws.Range("C10").Value = "School Without Walls @ Francis Stevens"
$ws.Range("D10").Value = "K-12 School"
$ws.Range("E10").Value = "2425 N STREET NW"
$ws.Range("H10").Value = 20037
$ws.Range("I10").Value = 1924
$ws.Range("J10").Value = "DISTRICT OF COLUMBIA"
$ws.Range("L10").Value = 127991
$ws.Range("M10").Value = 4.5
$ws.Range("N10").Value = 70.1
$ws.Range("P10").Value = 60

# ---------------------------------------------------------------------------
# Clear the leftover date format on the "Year Built" column so it shows
# plain numbers instead of dates (e.g. I2 was formatted as m/d/yyyy).
# ---------------------------------------------------------------------------
$ws.Range("I2:I10").Style = "Normal"

# Re-assert the Year Built values after clearing the style, in case the
# style reset also clears cell contents on this host.
$ws.Range("I2").Value = 1985
$ws.Range("I3").Value = 1991
$ws.Range("I4").Value = 1991
$ws.Range("I5").Value = 1962
$ws.Range("I6").Value = 1921
$ws.Range("I7").Value = 2004
$ws.Range("I8").Value = 1912
$ws.Range("I9").Value = 1880
$ws.Range("I10").Value = 1924

# ---------------------------------------------------------------------------
# Update the active selection to match the new view (M22).
# ---------------------------------------------------------------------------
$ws.Range("M22").Select() | Out-Null
